# Update L1cam-Itga5.xlsx with new TPM-derived values.
#
# For each data row the "Sending cluster" (col A) determines the Ligand
# average/total expression and its derived-specificity values (cols G,H,I,J)
# and the "Target cluster" (col D) determines the Receptor average/total
# expression and its derived-specificity values (cols M,N,O,P). The edge
# weight/specificity columns (Q,R,S,T) are simply the products of the
# matching ligand/receptor columns (Q=G*M, R=H*N, S=I*O, T=J*P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-"Sending cluster" ligand values: G, H, I, J
$sendingNew = @{
    "ECs"               = @(6.713252999999999, 20.139759, 0.3101840064655811, 0.3231642354899327)
    "FAPs"              = @(0.3343473333333333, 1.003042, 0.0154484264788496, 0.01609489473505086)
    "Inflammatory-Mac"  = @(6.661784666666667, 19.985354, 0.3078059262949933, 0.3206866401135023)
    "MuSCs"             = @(2.607918, 5.215835999999999, 0.1204981331366039, 0.08369373503331734)
    "Resolving-Mac"     = @(5.325505333333333, 15.976516, 0.2460635076239721, 0.2563604946281968)
}

# New per-"Target cluster" receptor values: M, N, O, P
$targetNew = @{
    "ECs"               = @(28.85518433333334, 86.56555300000001, 0.1999651185353207, 0.2044513327926365)
    "FAPs"              = @(51.17424933333334, 153.522748, 0.3546352265743414, 0.3625914622481308)
    "Inflammatory-Mac"  = @(29.393479, 88.180437, 0.2036954761578358, 0.2082653809291453)
    "MuSCs"             = @(9.499066500000001, 18.998133, 0.0658280999596015, 0.04486996822421697)
    "Resolving-Mac"     = @(25.37910966666666, 76.137329, 0.1758760787729007, 0.1798218558058706)
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2
    $target  = $ws.Cells.Item($r, 4).Value2

    if (-not $sending -or -not $sendingNew.ContainsKey($sending)) { continue }
    if (-not $target -or -not $targetNew.ContainsKey($target)) { continue }

    $g = $sendingNew[$sending][0]
    $h = $sendingNew[$sending][1]
    $i = $sendingNew[$sending][2]
    $j = $sendingNew[$sending][3]

    $m = $targetNew[$target][0]
    $n = $targetNew[$target][1]
    $o = $targetNew[$target][2]
    $p = $targetNew[$target][3]

    $ws.Cells.Item($r, 7).Value2  = $g   # G
    $ws.Cells.Item($r, 8).Value2  = $h   # H
    $ws.Cells.Item($r, 9).Value2  = $i   # I
    $ws.Cells.Item($r, 10).Value2 = $j   # J

    $ws.Cells.Item($r, 13).Value2 = $m   # M
    $ws.Cells.Item($r, 14).Value2 = $n   # N
    $ws.Cells.Item($r, 15).Value2 = $o   # O
    $ws.Cells.Item($r, 16).Value2 = $p   # P

    $ws.Cells.Item($r, 17).Value2 = ($g * $m)   # Q
    $ws.Cells.Item($r, 18).Value2 = ($h * $n)   # R
    $ws.Cells.Item($r, 19).Value2 = ($i * $o)   # S
    $ws.Cells.Item($r, 20).Value2 = ($j * $p)   # T
}
